# Edit: (1) change the table style on the Component-three table (slide 5,
#       shape 2) to the built-in style {9E0E6EC5-6B4C-4069-8E6A-F3BDF5AEBFAC};
#   (2) swap the deck's theme colours from the "Integral" (Red Violet)
#       scheme to the default "Office" scheme, mirroring the underlying
#       OOXML theme1.xml/theme2.xml content swap.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{9E0E6EC5-6B4C-4069-8E6A-F3BDF5AEBFAC}", $true)

# --- 2. Theme colours ------------------------------------------------------
# Office (default) theme colour values, in PpColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
